$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.238.86'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '3.278.43'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("D5").Value = "'587.50"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = '  +2.21%  '
$ws.Range("D6").Value = "'179.13"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = '  -1.97%  '
$ws.Range("E7").Value = '  +4.24%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("D10").Value = "'6.73"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = '  +1.54%  '
$ws.Range("D11").Value = "'0.402"
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = '3.850.58'
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("E13").Value = '  -3.43%  '
$ws.Range("D14").Value = '66.230.72'
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = "'26.34"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = '  -2.74%  '
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").Value = '3.285.36'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").Value = "'432.64"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").Value = "'5.52"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = '  -2.30%  '
$ws.Range("D20").Value = "'13.14"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = '  -2.96%  '
$ws.Range("D21").Value = "'7.39"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = '  -3.84%  '
$ws.Range("D22").Value = "'71.66"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = '  -2.97%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").Value = '3.422.60'
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  -4.46%  '
$ws.Range("D28").Value = "'8.83"
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = "'1.96"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  -2.68%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("E33").Value = '  -2.36%  '
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("E35").Value = '  -1.89%  '
$ws.Range("D36").Value = "'157.73"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = '  -2.19%  '
$ws.Range("E37").Value = '  -4.80%  '
$ws.Range("D38").Value = "'26.48"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = '  -3.63%  '
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("D40").Value = '2.778.66'
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("D41").Value = "'0.774"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = "'40.19"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").Value = "'6.03"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = '  -2.38%  '
$ws.Range("D45").Value = "'0.0659"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("D46").Value = "'320.59"
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").Value = "'2.29"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = "'23.18"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = '  -4.33%  '
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("E50").Value = '  +3.34%  '
$ws.Range("E51").Value = '  +0.09%  '
